$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aerobox")

# B1: was a time value (0.75 -> 18:00) formatted with the time style (s=3).
# Now it becomes literal text '"18:00"' and picks up the plain/general
# left-aligned style (same style as C1, i.e. s=2).
$ws.Range("B1").Value = '"18:00"'
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# B2: was a time value (0.75 -> 18:00). Now becomes the same literal text,
# but keeps its original (time-format) style (s=3).
$ws.Range("B2").Value = '"18:00"'

# B3: was a time value (0.41666... -> 10:00). Now becomes literal text,
# keeping its original style (s=3).
$ws.Range("B3").Value = '"10:00"'

# Column B is now wide enough to show the new text, with an explicit
# (best-fit) width.
$ws.Columns.Item(2).ColumnWidth = 11

# Selection moves to D7.
$ws.Range("D7").Select() | Out-Null
